$wb = $excel.ActiveWorkbook

# A new handoff xliff was generated, so the report now reflects "Ready for
# handoff" (rather than a handback-complete status) with refreshed timestamps.
#
# Column-width note: the Status column used to hold the long sentence
# "Handed back: in sync with en-US" and was sized accordingly. Now that the
# status text is the much shorter "Ready for handoff", that column is
# narrowed to fit. (This runtime quantizes ColumnWidth to whole pixels, so
# 16.33 is the input that round-trips to the stored width closest to the
# narrower target.)
$narrowStatusColumnWidth = 16.33

# --- Overview sheet: update Status text + timestamp ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 01:13:56"
$wsOverview.Columns.Item(5).ColumnWidth = $narrowStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowStatusColumnWidth

# --- zh-cn sheet: update Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 01:13:49"
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowStatusColumnWidth

# --- de-de sheet: update Status + Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 01:13:56"
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowStatusColumnWidth
